# Generate Report for Handback
#
# - "Status" column (C) on the zh-cn / de-de sheets moves from "In Translation"
#   to "Handed back: in sync with en-US" for both rows.
# - "Latest Target File" (I) and "Latest Handback File" (J) columns get filled
#   in for both rows/sheets, with I turning into a hyperlink to the source .md
#   (mirroring column A's hyperlink).
# - "Latest Handback DateTime" (K) is stamped for both rows/sheets.
# - Columns that now hold long file names / hyperlinks are widened.

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")
$ov = $wb.Worksheets.Item("Overview")

$zhUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2d98b6e44b3408105d8e5f8e6a40e828f11e4da/e2e/1d3b7a39-7304-47a2-843c-d1ad05ecf74b.md"
$zhUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2d98b6e44b3408105d8e5f8e6a40e828f11e4da/e2e/d17a0600-f08f-4403-ba57-a872c6db824f.md"

# --- Status column: "In Translation" -> "Handed back: in sync with en-US" ---
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$zh.Range("J2").Value = "1d3b7a39-7304-47a2-843c-d1ad05ecf74b.00a9c5d04915e7f6e4924d8e73b78dc19c305377.zh-cn.xlf"
$zh.Range("J3").Value = "d17a0600-f08f-4403-ba57-a872c6db824f.725eca774c16a80be44c811094c209f2f9fea53f.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-23 16:26:16"
$zh.Range("K3").Value = "2016-08-23 16:26:16"

$zh.Hyperlinks.Add($zh.Range("I2"), $zhUrl1, "", "", "1d3b7a39-7304-47a2-843c-d1ad05ecf74b.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), $zhUrl2, "", "", "d17a0600-f08f-4403-ba57-a872c6db824f.md") | Out-Null

# --- de-de: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$de.Range("J2").Value = "1d3b7a39-7304-47a2-843c-d1ad05ecf74b.00a9c5d04915e7f6e4924d8e73b78dc19c305377.de-de.xlf"
$de.Range("J3").Value = "d17a0600-f08f-4403-ba57-a872c6db824f.725eca774c16a80be44c811094c209f2f9fea53f.de-de.xlf"
$de.Range("K2").Value = "2016-08-23 16:26:24"
$de.Range("K3").Value = "2016-08-23 16:26:24"

$de.Hyperlinks.Add($de.Range("I2"), $zhUrl1, "", "", "1d3b7a39-7304-47a2-843c-d1ad05ecf74b.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), $zhUrl2, "", "", "d17a0600-f08f-4403-ba57-a872c6db824f.md") | Out-Null

# --- Column width adjustments (report columns now show longer file names) ---
$ov.Range("E1").ColumnWidth = 29.14
$ov.Range("F1").ColumnWidth = 29.14

$zh.Range("C1").ColumnWidth = 29.14
$zh.Range("I1").ColumnWidth = 39.14
$zh.Range("J1").ColumnWidth = 39.14

$de.Range("C1").ColumnWidth = 29.14
$de.Range("I1").ColumnWidth = 39.14
$de.Range("J1").ColumnWidth = 39.14
